# Add a "Date" column (new column B) to the Weekly Shop Schedule sheet,
# shifting the existing Time Slot / Technician / Status columns one to the
# right (old B/C/D -> new C/D/E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("B:B").Insert()

# Header for the new column, matching the bold style used by the other
# header cells in row 1.
$ws.Range("B1").Value = "Date"
$ws.Range("B1").Font.Bold = $true

# Fill in the dates for each day block (Monday - Friday of that week).
$ws.Range("B2:B4").Value = 45810
$ws.Range("B5:B7").Value = 45811
$ws.Range("B8:B10").Value = 45812
$ws.Range("B11:B13").Value = 45813
$ws.Range("B14:B16").Value = 45814

# Apply a date number format to B2, then copy that formatting across the
# rest of the column so every cell shares the same cell style.
$ws.Range("B2").NumberFormat = "mm-dd-yy"
$ws.Range("B2").Copy()
$ws.Range("B3:B16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F11").Select()
